$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stations")

# Re-run with updated catchments: point station fpath_or_id entries at the
# newer shapefiles for Samnanger_Storelva and Samnanger_Frolandskanalen.
$ws.Range("F4").Value = "/home/jovyan/projects/critical_loads_2/cl_vestland/shapefiles/Samnanger_Storelva_ny.shp"
$ws.Range("F6").Value = "/home/jovyan/projects/critical_loads_2/cl_vestland/shapefiles/Samnanger_Frolandskanalen_nyere.shp"

# Update the active selection on the stations sheet (last place the cursor
# was left after the edit).
$ws.Activate()
$ws.Range("F12").Select()
